$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 45867
$ws.Range("A25").NumberFormat = $ws.Range("A24").NumberFormat()
$ws.Range("B25").Value = 6370.85986328125
$ws.Range("C25").Value = 6409.259765625
$ws.Range("D25").Value = 6363.919921875
$ws.Range("E25").Value = 6405.6201171875
$ws.Range("F25").Value = 5076120000
$ws.Range("G25").Value = -0.0029594423887242
$ws.Range("H25").Value = -1
$ws.Range("I25").Value = 'kfa private wealth group llc ha million stake in vanguard sp etf nysearcavoo werba rubin papier wealth management ha million stock position in avantis international equity etf nysearcaavde werba rubin papier wealth management sell share of dimensional u targeted value etf nysearcadfat werba rubin papier wealth management cut stock holding in vanguard sp etf nysearcavoo mackenzie financial corp ha million stock holding in marsh mclennan company inc nysemmc foundation investment advisor llc acquires share of the goldman sachs group inc nysegs kfa private wealth group llc purchase share of the goldman sachs group inc nysegs werba rubin papier wealth management ha million position in vanguard value etf nysearcavtv mackenzie financial corp sell share of marsh mclennan company inc nysemmc werba rubin papier wealth management raise holding in exxon mobil corporation nysexom vanguard growth etf nysearcavug share acquired by werba rubin papier wealth management werba rubin papier wealth management buy share of vanguard growth etf nysearcavug tariff on south korea product threaten the kbeauty boom in the u tariff on south korea product threaten the kbeauty boom in the u tariff on south korea product threaten the kbeauty boom in the u tariff on south korea product threaten the kbeauty boom in the u tariff on south korea product threaten the kbeauty boom in the u how gas price have changed in nebraska in the last week how gas price have changed in alabama in the last week the jefferson to debut august th on fetv wall street coasting ahead of a week packed with potential flashpoint wall street coasting ahead of a week packed with potential flashpoint wall street coasting ahead of a week packed with potential flashpoint google launch free ai app builder opal that embrace vibe coding jim cramer explains why trump trade deal didnt bring on a market rally amphenol corporation nyseaph share purchased by fourthought financial partner llc td asset management inc lower holding in pool corporation nasdaqpool milton hydro power up outage notification with sm alert pilot dear microsoft stock fan mark your calendar for aug alibaba group holding limited nysebaba share sold by norden group llc norden group llc buy share of transdigm group incorporated nysetdg werba rubin papier wealth management ha stake in ishares california muni bond etf nysearcacmf contrasting quantum service nysepwr and mayville engineering nysemec law taking effect this week empowers minnesotan to request delete personal consumer data chinese export to u could decrease by billion by tariff simulator hedge fund and insider trading news ray dalio ken griffin michael gelband warren buffett de shaw lhc capital rli corp rli and more ethereum news today crypto market cap nears trillion a eth lead etp inflow nosaintstrainingcampjpg nosaintstrainingcampjpg nosaintstrainingcampjpg nosaintstrainingcampjpg noartgretnaadv copyjpg noartgretnaadv copyjpg tradebe environmental service partner with boy girl club of east chicago to support youth development and scholarship program noartgretnaadv copyjpg noartgretnaadv copyjpg noartgretnaadv copyjpg noartgretnaadv copyjpg mediolanum international fund ltd sell share of marriott international inc nasdaqmar mediolanum international fund ltd cut stock holding in edward lifesciences corporation nyseew joe rogan isnt letting go of epstein thats a problem for trump top sign your startup is built to die and how to fix it fast best ai tool to generate marketing content amazon is selling a convenient portable in steamer and iron for and shopper call it amazing how robot like robotic dog and magnetic slime are changing life ap business summarybrief at pm edt liven up your home this summer with up to off at kohl senator introduce bill to tighten aml regulation on u art market full pour magazine win best cocktail and spirit publication at tale of the cocktail annual spirited award purkiss capital advisor llc raise stock holding in analog device inc nasdaqadi baldwin investment management llc buy share of lowes company inc nyselow nyl investor llc ha million stake in vanguard total bond market etf nasdaqbnd morgan stanley nysems share sold by baldwin investment management llc purkiss capital advisor llc raise stock holding in philip morris international inc nysepm baldwin investment management llc increase stock holding in intercontinental exchange inc nyseice airbnb inc nasdaqabnb position reduced by norden group llc awm capital llc acquires share of vanguard real estate etf nysearcavnq baldwin investment management llc ha million stock holding in enterprise product partner lp nyseepd baldwin investment management llc grows position in intercontinental exchange inc nyseice uschina trade talk are showing sign of detente say mcc globals michelle carusocabrera candidate view duluth need to embrace growth in order to survive lazard freres gestion sa sell share of tmobile u inc nasdaqtmus solidarity wealth llc sell share of rtx corporation nysertx lazard freres gestion sa reduces stock position in procter gamble company the nysepg ishares tip bond etf nysearcatip share sold by awm capital llc xrp news today xrp price hold above a weekly close nears historic high lazard freres gestion sa ha million holding in becton dickinson and company nysebdx wincap financial llc ha position in qualcomm incorporated nasdaqqcom solidarity wealth llc decrease holding in nextera energy inc nysenee solidarity wealth llc cut stake in asml holding nv nasdaqasml lazard freres gestion sa cut stock position in parkerhannifin corporation nyseph avalon capital management ha million stock holding in spdr gold share nysearcagld avalon capital management boost stock position in spdr gold share nysearcagld solidarity wealth llc sell share of asml holding nv nasdaqasml samsungs huge gaming monitor is off for a limited time fiserv stock get pummeled after earnings look cheap here rayban maker essilorluxottica sale top estimate bloombergcom avalon capital management make new million investment in vanguard ftse europe etf nysearcavgk oxler private wealth llc make new investment in invesco nasdaq etf nasdaqqqqm invesco qqq nasdaqqqq stake increased by oxler private wealth llc oxler private wealth llc cut position in salesforce inc nysecrm evermay wealth management llc ha stake in oracle corporation nyseorcl solidarity wealth llc take position in novo nordisk a nysenvo oxler private wealth llc sell share of vanguard dividend appreciation etf nysearcavig church dwight co inc nysechd share sold by oxler private wealth llc evermay wealth management llc ha holding in cocacola company the nyseko share in vanguard intermediateterm treasury etf nasdaqvgit purchased by avalon capital management agp franklin llc increase stock position in the tjx company inc nysetjx att inc nyset share acquired by wincap financial llc yearold retail chain to close major location after historic downfall exxon mobil corporation nysexom share purchased by werba rubin papier wealth management kfa private wealth group llc ha stake in the goldman sachs group inc nysegs werba rubin papier wealth management ha million position in vanguard value etf nysearcavtv mackenzie financial corp ha million position in msci inc nysemsci pool corporation nasdaqpool share sold by td asset management inc td asset management inc grows stock position in martin marietta material inc nysemlm the jas gripen fighter jet engages in combat in thailand norden group llc acquires share of johnson johnson nysejnj caesar palace time square to fund hightech em ambucycletm program innovative lifesaving response model to dramatically reduce emergency response time set to come to hell kitchen and west midtown tradebe environmental service partner with boy girl club of east chicago to support youth development and scholarship program federal judge block big beautiful bill medicaid cut to planned parenthood ge vernova inc nysegev share acquired by hb wealth management llc joe rogan isnt letting go of epstein thats a problem for trump city of owasso investigating inappropriate police facebook post lexington police make second arrest in federal street homicide josh hawley team with democrat on legislation targeting ai and drug cost dimensional international core equity etf batsdfic is awm capital llcs th largest position ishares gold trust nysearcaiau stock position increased by purkiss capital advisor llc awm capital llc reduces holding in vanguard total bond market etf nasdaqbnd baldwin investment management llc buy share of enterprise product partner lp nyseepd baldwin investment management llc acquires share of intuitive surgical inc nasdaqisrg lindsey graham get republican challenger project author to announce senate bid will the boeing strike delay the f fighter jet lazard freres gestion sa sell share of parkerhannifin corporation nyseph wincap financial llc acquires share of qualcomm incorporated nasdaqqcom ishares core u aggregate bond etf nysearcaagg share acquired by foundation investment advisor llc live pm et unleash prosperity project hold panel on finding of it vote with your foot project bob gallagher littmann show he the right choice to lead annapolis commentary comparing taylor morrison home nysetmhc taylor wimpey otcmktstwody avalon capital management sell share of invesco qqq nasdaqqqq obama called out for shady gaza post failing to mention why innocent people are dying of starvation why scotland protest about trump should concern every american opinion obamaappointed judge order trump to continue funding planned parenthood congress eye shortterm govt funding stopgap a appropriation bill stall meridian management co sell share of visa inc nysev sunpointe llc decrease stake in accenture plc nyseacn eqt and cpp investment to acquire neogov a provider of hr and compliance software for u public sector agency from warburg pincus and carlyle agp franklin llc boost stock holding in palantir technology inc nasdaqpltr physician financial service inc boost stock holding in chevron corporation nysecvx schubert co take position in ge aerospace nysege share in ishares msci emerging market etf nysearcaeem purchased by avalon capital management schubert co sell share of cigna group nyseci trump call for more gaza food aid and talk up his golf course in meeting with uk prime minister empirical finance llc grows stake in ametek inc nyseame durbin demand recording of ghislaine maxwelldepartment of justice talk durbin demand recording of ghislaine maxwelldepartment of justice talk durbin demand recording of ghislaine maxwelldepartment of justice talk durbin demand recording of ghislaine maxwelldepartment of justice talk durbin demand recording of ghislaine maxwelldepartment of justice talk durbin demand recording of ghislaine maxwelldepartment of justice talk prosperity consulting group llc trim holding in oneok inc nyseoke jpmorgan ultrashort income etf nysearcajpst share sold by prosperity consulting group llc trump give russia day deadline to end ukraine democraticled state sue over trump administration effort to obtain personal information of snap recipient we are not dog palestinian condemn humiliation of chasing after airdropped aid in gaza ginnie graham dont let ryan walter pornogate scandal deflect from his other failing vantage financial partner llc purchase share of ishares u treasury bond etf batsgovt vantage financial partner llc sell share of ishares core sp u value etf nasdaqiusv rosen a leading investor right law firm encourages replimune group inc investor to secure counsel before important deadline in security class action first filed by the firm repl in supreme court appeal ghislaine maxwell fight justice department view of epstein plea deal cnn confronts dem senator with party dismal favorability and here what his response really mean sheaff brock investment advisor llc ha million holding in cummins inc nysecmi meridian management co decrease stock position in ishares russell etf nysearcaiwm ishares russell etf nysearcaiwm share purchased by portland global advisor llc pharmaceutical stock worth watching july th share in vanguard information technology etf nysearcavgt acquired by passive capital management llc space stock to watch today july th thailand and cambodia truce agreement come into effect paradigm strategy in wealth management llc buy share of jpmorgan nasdaq equity premium income etf nasdaqjepq new republic capital llc ha stock holding in tesla inc nasdaqtsla commissioner celebrate joe hewgleys legacy of resilience and impact over last year a nonexhaustive timeline of robert f kennedy jr war on vaccine ap business summarybrief at pm edt maine police chief want to unseat chellie pingree nkcfo llc ha million stake in expedia group inc nasdaqexpe nkcfo llc cut stock holding in agnico eagle mine limited nyseaem nkcfo llc purchase new stake in range resource corporation nyserrc trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story trump seek quick murdoch deposition in wall street journal lawsuit over epstein story dead target u terrorgram leader linked to white supremacist attack to plead guilty letter is america backsliding on liberty equality letter is america backsliding on liberty equality ghislaine maxwell pleads with supreme court trump to intervene in criminal case ghislaine maxwell pleads with supreme court trump to intervene in criminal case ghislaine maxwell asks supreme court to toss out sextrafficking conviction ghislaine maxwell pleads with supreme court trump to intervene in criminal case ghislaine maxwell pleads with supreme court trump to intervene in criminal case'
$ws.Range("J25").Value = 'kfa private wealth million stake vanguard sp nysearcavoo werba rubin papier wealth million stock avantis international equity nysearcaavde werba rubin papier wealth sell dimensional targeted value nysearcadfat werba rubin papier wealth cut stock vanguard sp nysearcavoo mackenzie financial corp million stock marsh mclennan company nysemmc foundation investment advisor acquires goldman sachs nysegs kfa private wealth purchase goldman sachs nysegs werba rubin papier wealth million vanguard value nysearcavtv mackenzie financial corp sell marsh mclennan company nysemmc werba rubin papier wealth raise exxon mobil corporation nysexom vanguard growth nysearcavug acquired werba rubin papier wealth werba rubin papier wealth buy vanguard growth nysearcavug tariff south korea product threaten kbeauty boom tariff south korea product threaten kbeauty boom tariff south korea product threaten kbeauty boom tariff south korea product threaten kbeauty boom tariff south korea product threaten kbeauty boom gas price changed nebraska gas price changed alabama jefferson debut august fetv wall street coasting ahead packed potential flashpoint wall street coasting ahead packed potential flashpoint wall street coasting ahead packed potential flashpoint google launch free ai app builder opal embrace vibe coding jim cramer explains trump trade deal didnt market rally amphenol corporation nyseaph purchased fourthought financial partner td asset lower pool corporation nasdaqpool milton hydro power outage notification sm alert pilot dear microsoft stock fan mark calendar aug alibaba limited nysebaba norden norden buy transdigm incorporated nysetdg werba rubin papier wealth stake ishares california muni bond nysearcacmf contrasting quantum service nysepwr mayville engineering nysemec taking effect empowers minnesotan request delete personal consumer data chinese export decrease billion tariff simulator hedge fund insider trading ray dalio ken griffin michael gelband warren buffett de shaw lhc rli corp rli ethereum crypto market cap nears trillion eth lead etp inflow nosaintstrainingcampjpg nosaintstrainingcampjpg nosaintstrainingcampjpg nosaintstrainingcampjpg noartgretnaadv copyjpg noartgretnaadv copyjpg tradebe environmental service partner boy girl club east chicago support youth development scholarship program noartgretnaadv copyjpg noartgretnaadv copyjpg noartgretnaadv copyjpg noartgretnaadv copyjpg mediolanum international fund ltd sell marriott international nasdaqmar mediolanum international fund ltd cut stock edward lifesciences corporation nyseew joe rogan isnt letting epstein thats problem trump sign startup built die fix fast best ai tool generate marketing content amazon selling convenient portable steamer iron shopper call amazing robot robotic dog magnetic slime changing life business pm liven summer kohl senator introduce bill tighten aml regulation art market full pour magazine best cocktail spirit publication tale cocktail annual spirited award purkiss advisor raise stock analog device nasdaqadi baldwin investment buy lowes company nyselow nyl investor million stake vanguard total bond market nasdaqbnd morgan stanley nysems baldwin investment purkiss advisor raise stock philip morris international nysepm baldwin investment increase stock intercontinental exchange nyseice airbnb nasdaqabnb reduced norden awm acquires vanguard real estate nysearcavnq baldwin investment million stock enterprise product partner lp nyseepd baldwin investment grows intercontinental exchange nyseice uschina trade showing sign detente mcc globals michelle carusocabrera candidate view duluth embrace growth order survive lazard freres gestion sa sell tmobile nasdaqtmus solidarity wealth sell rtx corporation nysertx lazard freres gestion sa reduces stock procter gamble company nysepg ishares tip bond nysearcatip awm xrp xrp price weekly close nears historic high lazard freres gestion sa million becton dickinson company nysebdx wincap financial qualcomm incorporated nasdaqqcom solidarity wealth decrease nextera energy nysenee solidarity wealth cut stake asml nv nasdaqasml lazard freres gestion sa cut stock parkerhannifin corporation nyseph avalon million stock spdr gold nysearcagld avalon boost stock spdr gold nysearcagld solidarity wealth sell asml nv nasdaqasml samsungs gaming monitor limited fiserv stock pummeled earnings look cheap rayban maker essilorluxottica sale estimate bloombergcom avalon million investment vanguard ftse europe nysearcavgk oxler private wealth investment invesco nasdaq nasdaqqqqm invesco qqq nasdaqqqq stake increased oxler private wealth oxler private wealth cut salesforce nysecrm evermay wealth stake oracle corporation nyseorcl solidarity wealth take novo nordisk nysenvo oxler private wealth sell vanguard dividend appreciation nysearcavig church dwight co nysechd oxler private wealth evermay wealth cocacola company nyseko vanguard intermediateterm treasury nasdaqvgit purchased avalon agp franklin increase stock tjx company nysetjx att nyset acquired wincap financial yearold retail chain close major location historic downfall exxon mobil corporation nysexom purchased werba rubin papier wealth kfa private wealth stake goldman sachs nysegs werba rubin papier wealth million vanguard value nysearcavtv mackenzie financial corp million msci nysemsci pool corporation nasdaqpool td asset td asset grows stock martin marietta material nysemlm jas gripen fighter jet engages combat thailand norden acquires johnson johnson nysejnj caesar palace square fund hightech em ambucycletm program innovative lifesaving response model dramatically reduce emergency response set come hell kitchen west midtown tradebe environmental service partner boy girl club east chicago support youth development scholarship program federal judge block bill medicaid cut planned parenthood ge vernova nysegev acquired hb wealth joe rogan isnt letting epstein thats problem trump owasso investigating inappropriate police facebook post lexington police arrest federal street homicide josh hawley team democrat legislation targeting ai drug cost dimensional international core equity batsdfic awm llcs largest ishares gold trust nysearcaiau stock increased purkiss advisor awm reduces vanguard total bond market nasdaqbnd baldwin investment buy enterprise product partner lp nyseepd baldwin investment acquires intuitive surgical nasdaqisrg lindsey graham republican challenger project author announce senate bid boeing strike delay f fighter jet lazard freres gestion sa sell parkerhannifin corporation nyseph wincap financial acquires qualcomm incorporated nasdaqqcom ishares core aggregate bond nysearcaagg acquired foundation investment advisor pm et unleash prosperity project panel finding vote foot project bob gallagher littmann choice lead annapolis commentary comparing taylor morrison nysetmhc taylor wimpey otcmktstwody avalon sell invesco qqq nasdaqqqq obama shady gaza post failing mention innocent people dying starvation scotland protest trump concern american opinion obamaappointed judge order trump continue funding planned parenthood congress eye shortterm govt funding stopgap appropriation bill stall meridian co sell visa nysev sunpointe decrease stake accenture plc nyseacn eqt cpp investment acquire neogov provider hr compliance software public sector agency warburg pincus carlyle agp franklin boost stock palantir technology nasdaqpltr physician financial service boost stock chevron corporation nysecvx schubert co take ge aerospace nysege ishares msci emerging market nysearcaeem purchased avalon schubert co sell cigna nyseci trump call gaza aid golf course meeting uk prime minister empirical finance grows stake ametek nyseame durbin demand recording ghislaine maxwelldepartment justice durbin demand recording ghislaine maxwelldepartment justice durbin demand recording ghislaine maxwelldepartment justice durbin demand recording ghislaine maxwelldepartment justice durbin demand recording ghislaine maxwelldepartment justice durbin demand recording ghislaine maxwelldepartment justice prosperity consulting trim oneok nyseoke jpmorgan ultrashort income nysearcajpst prosperity consulting trump give russia deadline end ukraine democraticled sue trump administration effort obtain personal information snap recipient dog palestinian condemn humiliation chasing airdropped aid gaza ginnie graham dont ryan walter pornogate scandal deflect failing vantage financial partner purchase ishares treasury bond batsgovt vantage financial partner sell ishares core sp value nasdaqiusv rosen leading investor firm encourages replimune investor secure counsel important deadline security class action filed firm repl supreme appeal ghislaine maxwell fight justice department view epstein plea deal cnn confronts dem senator party dismal favorability response sheaff brock investment advisor million cummins nysecmi meridian co decrease stock ishares russell nysearcaiwm ishares russell nysearcaiwm purchased portland global advisor pharmaceutical stock worth watching vanguard information technology nysearcavgt acquired passive space stock watch thailand cambodia truce agreement come effect paradigm strategy wealth buy jpmorgan nasdaq equity premium income nasdaqjepq republic stock tesla nasdaqtsla commissioner celebrate joe hewgleys legacy resilience impact nonexhaustive timeline robert f kennedy jr war vaccine business pm maine police chief unseat chellie pingree nkcfo million stake expedia nasdaqexpe nkcfo cut stock agnico eagle mine limited nyseaem nkcfo purchase stake range resource corporation nyserrc trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story trump seek quick murdoch deposition wall street journal lawsuit epstein story dead target terrorgram linked white supremacist attack plead guilty america backsliding liberty equality america backsliding liberty equality ghislaine maxwell pleads supreme trump intervene criminal ghislaine maxwell pleads supreme trump intervene criminal ghislaine maxwell asks supreme toss sextrafficking conviction ghislaine maxwell pleads supreme trump intervene criminal ghislaine maxwell pleads supreme trump intervene criminal'
